$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "60.969.73"
Set-TextCell "E2" "  +0.26%  "

Set-TextCell "D3" "3.374.08"
Set-TextCell "E3" "  -0.35%  "

Set-TextCell "D4" "1.00"
Set-TextCell "E4" "  +0.07%  "

Set-TextCell "D5" "572.86"
Set-TextCell "E5" "  +0.69%  "

Set-TextCell "D6" "138.76"
Set-TextCell "E6" "  -1.61%  "

Set-TextCell "E7" "  -0.01%  "

Set-TextCell "E8" "  -0.37%  "

Set-TextCell "D9" "7.68"
Set-TextCell "E9" "  +2.54%  "

Set-TextCell "E10" "  -1.68%  "

Set-TextCell "D11" "0.383"
Set-TextCell "E11" "  -2.83%  "

Set-TextCell "D12" "3.949.05"
Set-TextCell "E12" "  -0.37%  "

Set-TextCell "E13" "  +0.81%  "

Set-TextCell "D14" "28.12"
Set-TextCell "E14" "  -1.24%  "

Set-TextCell "D15" "3.366.27"
Set-TextCell "E15" "  -0.57%  "

Set-TextCell "E16" "  -1.46%  "

Set-TextCell "D17" "61.034.23"
Set-TextCell "E17" "  +0.27%  "

Set-TextCell "D18" "6.10"
Set-TextCell "E18" "  -1.67%  "

Set-TextCell "D19" "13.53"
Set-TextCell "E19" "  -3.15%  "

Set-TextCell "D20" "8.92"
Set-TextCell "E20" "  -0.72%  "

Set-TextCell "D21" "385.32"
Set-TextCell "E21" "  +0.47%  "

Set-TextCell "D22" "74.97"
Set-TextCell "E22" "  +1.56%  "

Set-TextCell "E23" "  -1.29%  "

Set-TextCell "E24" "  +0.25%  "

Set-TextCell "D25" "0.0000111"
Set-TextCell "E25" "  -5.00%  "

Set-TextCell "E26" "  +5.72%  "

Set-TextCell "E27" "  +0.10%  "

Set-TextCell "D28" "7.14"
Set-TextCell "E28" "  -3.46%  "

Set-TextCell "D29" "7.93"
Set-TextCell "E29" "  -0.57%  "

Set-TextCell "D30" "2.12"
Set-TextCell "E30" "  -1.17%  "

Set-TextCell "E31" "  -0.03%  "

Set-TextCell "D32" "1.34"
Set-TextCell "E32" "  -6.51%  "

Set-TextCell "D33" "23.04"
Set-TextCell "E33" "  -2.18%  "

Set-TextCell "D34" "6.84"
Set-TextCell "E34" "  -1.81%  "

Set-TextCell "D35" "167.05"
Set-TextCell "E35" "  +1.06%  "

Set-TextCell "D36" "4.93"
Set-TextCell "E36" "  -0.80%  "

Set-TextCell "D37" "3.409.70"
Set-TextCell "E37" "  -0.19%  "

Set-TextCell "E38" "  -2.54%  "

Set-TextCell "E39" "  -1.92%  "

Set-TextCell "D40" "25.71"
Set-TextCell "E40" "  -8.79%  "

Set-TextCell "D41" "0.773"
Set-TextCell "E41" "  -0.61%  "

Set-TextCell "D42" "4.35"
Set-TextCell "E42" "  -1.31%  "

Set-TextCell "D43" "1.64"
Set-TextCell "E43" "  -1.36%  "

Set-TextCell "E44" "  -1.02%  "

Set-TextCell "D45" "2.452.36"
Set-TextCell "E45" "  -1.40%  "

Set-TextCell "E46" "  -2.35%  "

Set-TextCell "D47" "0.999"
Set-TextCell "E47" "  -0.03%  "

Set-TextCell "D48" "22.08"
Set-TextCell "E48" "  -6.07%  "

Set-TextCell "E49" "  -4.54%  "

Set-TextCell "D50" "2.00"
Set-TextCell "E50" "  -3.32%  "

Set-TextCell "E51" "  -2.86%  "
